# Updated cryptos list on Sat Aug 26 05:50:15 UTC 2023 with GitHub Actions
#
# Applies the per-row Price (D) / Volume(1h) (E) refresh, plus the two
# row-pair swaps (WrappedEther <-> Polkadot at rows 12/13, and
# EnergySwap <-> Cronos at rows 48/49) captured in the diff.
#
# Quirks of this COM-interop engine worked around below:
#  * Function parameters bind POSITIONALLY ONLY - named args (-Row, ...)
#    are silently ignored/blank, so every helper call below uses plain
#    positional args, with an explicit $null standing in for "leave this
#    cell alone".
#  * Assigning a plain numeric-looking string (e.g. "217.62") to
#    Range.Value auto-coerces it to a real number, same as typing it into
#    Excel - but the source sheet stores these Price cells as TEXT
#    (t="inlineStr"), even the ones that happen to look like plain
#    decimals. Forcing NumberFormat="@" immediately before the write (and
#    resetting the style back to "Normal" immediately after, so no stray
#    number-format style sticks to the cell) keeps the write text-typed
#    without changing the cell's on-screen style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText {
    param([int]$Row, [string]$Val)
    $c = $ws.Cells.Item($Row, 4)
    $c.NumberFormat = "@"
    $c.Value = $Val
    $c.Style = "Normal"
}

function Set-Row {
    param(
        [int]$Row,
        $Coin,
        $Link,
        $Price,
        $Volume
    )

    if ($Coin -ne $null) { $ws.Cells.Item($Row, 2).Value = $Coin }
    if ($Link -ne $null) { $ws.Cells.Item($Row, 3).Value = $Link }
    if ($Price -ne $null) { Set-PriceText $Row $Price }
    if ($Volume -ne $null) { $ws.Cells.Item($Row, 5).Value = $Volume }
}

# Row 2 - Bitcoin
Set-Row 2 $null $null "26.099.70" "  -0.25%  "
# Row 3 - Ethereum
Set-Row 3 $null $null "1.654.77" "  -0.31%  "
# Row 5 - BNB
Set-Row 5 $null $null "217.62" "  +0.59%  "
# Row 6 - XRP
Set-Row 6 $null $null "0.5238" "  +0.64%  "
# Row 7 - USDC
Set-Row 7 $null $null "1.002" "  -0.23%  "
# Row 8 - Cardano (price unchanged)
Set-Row 8 $null $null $null "  -1.18%  "
# Row 9 - Dogecoin
Set-Row 9 $null $null "0.06342" "  +1.08%  "
# Row 10 - Solana
Set-Row 10 $null $null "20.33" "  -2.19%  "
# Row 11 - TRON
Set-Row 11 $null $null "0.07800" "  +1.08%  "

# Rows 12/13 swap: WrappedEther <-> Polkadot
Set-Row 12 "Polkadot" "https://coinranking.com/coin/25W7FG7om+polkadot-dot" "4.498" "  +1.59%  "
Set-Row 13 "WrappedEther" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth" "1.659.44" "  +0.38%  "

# Row 14 - Polygon
Set-Row 14 $null $null "0.5466" "  +0.55%  "
# Row 15 - ShibaInu
Set-Row 15 $null $null "0.0₅8190" "  +0.28%  "
# Row 16 - Litecoin (price unchanged)
Set-Row 16 $null $null $null "  +1.23%  "
# Row 17 - WrappedBTC
Set-Row 17 $null $null "26.109.15" "  -0.34%  "
# Row 18 - Dai (price unchanged)
Set-Row 18 $null $null $null "  -0.30%  "
# Row 19 - Uniswap
Set-Row 19 $null $null "4.574" "  -1.10%  "
# Row 20 - BitcoinCash
Set-Row 20 $null $null "191.17" "  -0.42%  "
# Row 21 - Avalanche
Set-Row 21 $null $null "10.04" "  -0.27%  "
# Row 22 - Chainlink
Set-Row 22 $null $null "6.029" "  -0.52%  "
# Row 23 - BinanceUSD (price unchanged)
Set-Row 23 $null $null $null "  -0.29%  "
# Row 24 - Monero
Set-Row 24 $null $null "142.16" "  +1.67%  "
# Row 25 - Stellar
Set-Row 25 $null $null "0.1239" "  +0.79%  "
# Row 26 - Cosmos
Set-Row 26 $null $null "7.229" "  +0.66%  "
# Row 27 - EthereumClassic (price unchanged)
Set-Row 27 $null $null $null "  +0.56%  "
# Row 28 - Toncoin (price unchanged)
Set-Row 28 $null $null $null "  +1.52%  "
# Row 29 - Hedera
Set-Row 29 $null $null "0.05903" "  -1.48%  "
# Row 30 - PancakeSwap
Set-Row 30 $null $null "1.278" "  +0.61%  "
# Row 31 - InternetComputer(DFINITY)
Set-Row 31 $null $null "3.519" "  -0.87%  "
# Row 32 - Filecoin
Set-Row 32 $null $null "3.246" "  -0.15%  "
# Row 33 - LidoDAOToken (price unchanged)
Set-Row 33 $null $null $null "  -1.42%  "
# Row 34 - ARBITRUM
Set-Row 34 $null $null "0.9511" "  -1.71%  "
# Row 35 - MXToken (price unchanged)
Set-Row 35 $null $null $null "  +0.68%  "
# Row 36 - HuobiToken (price unchanged)
Set-Row 36 $null $null $null "  -0.30%  "
# Row 37 - ImmutableX
Set-Row 37 $null $null "0.5682" "  +0.20%  "
# Row 38 - VeChain
Set-Row 38 $null $null "0.01615" "  +1.04%  "
# Row 39 - FraxShare
Set-Row 39 $null $null "5.827" "  -3.04%  "
# Row 40 - TrustWalletToken
Set-Row 40 $null $null "0.8492" "  -0.78%  "
# Row 41 - PaxDollar (price unchanged)
Set-Row 41 $null $null $null "  -0.14%  "
# Row 42 - Maker
Set-Row 42 $null $null "1.029.65" "  +1.54%  "
# Row 43 - Quant
Set-Row 43 $null $null "102.68" "  +2.16%  "
# Row 44 - RocketPoolETH
Set-Row 44 $null $null "1.800.25" "  -0.03%  "
# Row 45 - Aave
Set-Row 45 $null $null "57.12" "  +0.52%  "
# Row 46 - Frax
Set-Row 46 $null $null "1.005" "  -0.39%  "
# Row 47 - Mantle
Set-Row 47 $null $null "0.4306" "  +2.57%  "

# Rows 48/49 swap: EnergySwap <-> Cronos
Set-Row 48 "Cronos" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro" "0.05165" "  -0.16%  "
Set-Row 49 "EnergySwap" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" "7.861" "  -1.82%  "

# Row 50 - RenderToken (price unchanged)
Set-Row 50 $null $null $null "  +1.35%  "
# Row 51 - Algorand
Set-Row 51 $null $null "0.09686" "  -0.39%  "
